# Updates the cryptos worksheet price/volume columns (D, E) to the
# latest GitHub Actions scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper writes Price values as literal text (e.g. "1.000", "0.9998")
# so trailing/insignificant-looking zeros survive. Cells whose new value would
# otherwise be auto-parsed as a number by Excel are pre-formatted as Text.
$textPriceRows = @(4, 5, 6, 7, 8, 9, 10, 11, 13, 14, 15, 18, 19, 22, 23, 24, 25, 27, 28, 29, 30, 31, 32, 33, 34, 35, 37, 38, 39, 42, 43, 44, 45, 47, 48, 50)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "29.190.82"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.860.58"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "0.7081"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "240.35"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "0.3078"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "0.07652"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").Value = "24.76"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "0.08426"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "1.858.33"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "5.183"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("D14").Value = "0.7093"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "91.09"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "29.210.03"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "243.00"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "0.000007819"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "2.111.60"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "7.855"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "8.913"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "18.43"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "1.317"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").Value = "4.402"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "4.216"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").Value = "0.05127"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("D34").Value = "0.8143"
$ws.Range("E34").Value = "  +12.67%  "
$ws.Range("D35").Value = "1.911"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("D37").Value = "2.676"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "0.01843"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").Value = "2.700"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").Value = "1.171.06"
$ws.Range("E40").Value = "  -6.74%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "0.8959"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "72.80"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "101.86"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "2.008.96"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "0.5171"
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("D48").Value = "1.772"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Value = "9.257"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("E51").Value = "  +0.30%  "
